$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168. This pushes the existing rows 168-275
# down to 169-276, preserving all of their data/formatting.
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new record.
$ws.Cells.Item(168, 1).Value = 5
$ws.Cells.Item(168, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(168, 3).Value = "Maule"
$ws.Cells.Item(168, 4).Value = Get-Date -Year 2022 -Month 2 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(168, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(168, 5).Value = 7
$ws.Cells.Item(168, 6).Value = 100114013
$ws.Cells.Item(168, 7).Value = "Zanahoria"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 300
$ws.Cells.Item(168, 11).Value = 8000
$ws.Cells.Item(168, 12).Value = 8000
$ws.Cells.Item(168, 13).Value = 8000
$ws.Cells.Item(168, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(168, 15).Value = "Región de Ñuble"
$ws.Cells.Item(168, 16).Value = 400
$ws.Cells.Item(168, 17).Value = 20
$ws.Cells.Item(168, 18).Value = "Hortaliza"
